$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 109, shifting existing rows
# 109-123 down to 110-124 (mirrors the xlsx diff, which effectively
# pushes the whole Papaya block down by one row and appends a fresh
# weekly entry at the top of the block).
$ws.Rows(109).Insert()

# Populate the new row 109 with this week's reading (columns A, B, C,
# E, F, G, H, I, J, L, R stay identical to the rest of the block).
$ws.Cells.Item(109, 1).Value = 10
$ws.Cells.Item(109, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(109, 3).Value = "La Araucanía"
$ws.Cells.Item(109, 4).Value = 45194
$ws.Cells.Item(109, 5).Value = 9
$ws.Cells.Item(109, 6).Value = "Fruta"
$ws.Cells.Item(109, 7).Value = 100108
$ws.Cells.Item(109, 8).Value = "Tropicales y subtropicales"
$ws.Cells.Item(109, 9).Value = 100108004
$ws.Cells.Item(109, 10).Value = "Papaya"
$ws.Cells.Item(109, 11).Value = "Cultivar IV Región"
$ws.Cells.Item(109, 12).Value = "Primera"
$ws.Cells.Item(109, 13).Value = 55
$ws.Cells.Item(109, 14).Value = 24000
$ws.Cells.Item(109, 15).Value = 24000
$ws.Cells.Item(109, 16).Value = 24000
$ws.Cells.Item(109, 17).Value = "$/bandeja 10 kilos"
$ws.Cells.Item(109, 18).Value = "Provincia del Elquí"
$ws.Cells.Item(109, 19).Value = 2400
$ws.Cells.Item(109, 20).Value = 10
